$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("gamma_frailty_fit")

$ws.Range("A2:A210").NumberFormat = "@"

$ws.Range("A2:A31").Value = "2022-07-25"
$ws.Range("A32:A61").Value = "2023-02-27"
$ws.Range("A62:A101").Value = "2023-08-28"
$ws.Range("A102:A151").Value = "2024-01-01"
$ws.Range("A152:A210").Value = "2024-06-24"
